$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header labels
$ws.Range("C1").Value = "rules"
$ws.Range("E1").Value = "adaptive_filter"

# Update column E (adaptive_filter) values from numeric 1 to string "RLS"
# and refresh precision of F, G, H (RMSE, NDEI, MAE) columns for rows 2-12

$ws.Range("E2").Value = "RLS"
$ws.Range("F2").Value = 0.606666490182016
$ws.Range("G2").Value = 1.693714455756965
$ws.Range("H2").Value = 0.4661378836489616

$ws.Range("E3").Value = "RLS"
$ws.Range("F3").Value = 0.5444379647031458
$ws.Range("G3").Value = 1.519982504397032
$ws.Range("H3").Value = 0.4470881296805422

$ws.Range("E4").Value = "RLS"
$ws.Range("F4").Value = 0.5600510251094211
$ws.Range("G4").Value = 1.563571637036913
$ws.Range("H4").Value = 0.4747365785701427

$ws.Range("E5").Value = "RLS"
$ws.Range("F5").Value = 0.6181671913588473
$ws.Range("G5").Value = 1.725822548341238
$ws.Range("H5").Value = 0.5274575743913941

$ws.Range("E6").Value = "RLS"
$ws.Range("F6").Value = 0.6757847624522083
$ws.Range("G6").Value = 1.886681462828424
$ws.Range("H6").Value = 0.5904330219172649

$ws.Range("E7").Value = "RLS"
$ws.Range("F7").Value = 0.6903711460265074
$ws.Range("G7").Value = 1.927404280252531
$ws.Range("H7").Value = 0.6107903267970767

$ws.Range("E8").Value = "RLS"
$ws.Range("F8").Value = 0.6438332028301864
$ws.Range("G8").Value = 1.797477890039094
$ws.Range("H8").Value = 0.5757396541367004

$ws.Range("E9").Value = "RLS"
$ws.Range("F9").Value = 0.5891664481288252
$ws.Range("G9").Value = 1.644857176376079
$ws.Range("H9").Value = 0.5298847891860214

$ws.Range("E10").Value = "RLS"
$ws.Range("F10").Value = 0.4900735525211847
$ws.Range("G10").Value = 1.368205882016436
$ws.Range("H10").Value = 0.4377291768399727

$ws.Range("E11").Value = "RLS"
$ws.Range("F11").Value = 0.3508336757110992
$ws.Range("G11").Value = 0.979470726889762
$ws.Range("H11").Value = 0.2856373854968268

$ws.Range("E12").Value = "RLS"
$ws.Range("F12").Value = 0.338942622392419
$ws.Range("G12").Value = 0.9462728344299641
$ws.Range("H12").Value = 0.2662458878999651
